$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item("Table 6")
$tbl = $shp.Table

# Row 4 ("PPT Link") -> point it at the same repo file the "Word Document
# Report Link" row already references, matching the author's re-upload.
$cell = $tbl.Cell(4, 2)
$cell.Shape.TextFrame.TextRange.Text = "https://github.com/dharrini06/WeatherPredictionAI/blob/main/Weather%20Prediction%20AI%20project.docx"

# Let PowerPoint's table auto-layout recompute the frame's rendered height
# now that row 4 wraps to fewer lines than before.
$shp.Height = $shp.Height
